$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from H1 (the last existing header cell) onto the two new
# header cells so I1/J1 pick up the same bold/centered/bordered format.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
